$wb = $excel.ActiveWorkbook
$occ = $wb.Worksheets.Item("Occupancy")
$cpue = $wb.Worksheets.Item("CPUE")

# --- Update the Longnose Dace model-selection table on "Occupancy" ---
# Row 9: top model now uses Area_km2 instead of avwid
$occ.Range("A9").Value = "p(~1)Psi(~Area_km2 + pctcbbl + pctSlope + med_len + BRT_100m)"
$occ.Range("B9").Value = 7
$occ.Range("C9").Value = 246.70421999999999
$occ.Range("D9").Value = 261.55918183206097
$occ.Range("E9").Value = 0
$occ.Range("F9").Value = 0.92818034409908401

# Row 10: now the med_len + BRT_100m model
$occ.Range("A10").Value = "p(~1)Psi(~med_len + BRT_100m)"
$occ.Range("B10").Value = 4
$occ.Range("C10").Value = 258.67410000000001
$occ.Range("D10").Value = 266.97260746268699
$occ.Range("E10").Value = 5.4134256306255697
$occ.Range("F10").Value = 0.061961500649566799

# Row 11: now the Area_km2 + pctcbbl + pctSlope model
$occ.Range("A11").Value = "p(~1)Psi(~Area_km2 + pctcbbl + pctSlope)"
$occ.Range("B11").Value = 5
$occ.Range("C11").Value = 260.54451
$occ.Range("D11").Value = 270.99563781954902
$occ.Range("E11").Value = 9.43645598748782
$occ.Range("F11").Value = 0.0082895696551185706

# Row 12: null model, same text, updated DeltaAICc/weight
$occ.Range("A12").Value = "p(~1)Psi(~1)"
$occ.Range("B12").Value = 2
$occ.Range("C12").Value = 270.23705000000001
$occ.Range("D12").Value = 274.32528529411798
$occ.Range("E12").Value = 12.7661034620566
$occ.Range("F12").Value = 0.0015685855962307999

# --- Update sheet selections / active tab: "Occupancy" becomes the active tab ---
$cpue.Range("A24").Select()
$occ.Activate()
$occ.Range("C9:F12").Select()
